# Update "Datos actualizados" timestamp in the title cell (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "Datos actualizados a 8 de Julio de 2020 a las 14:02"

# Row 6 - India
$ws.Cells.Item(6,2).Value = 746506
$ws.Cells.Item(6,3).Value = 3025
$ws.Cells.Item(6,4).Value = 458618
$ws.Cells.Item(6,5).Value = 267204

# Row 19 - Alemania
$ws.Cells.Item(19,2).Value = 198399
$ws.Cells.Item(19,3).Value = 44
$ws.Cells.Item(19,5).Value = 6596

# Row 34 - Emiratos Arabes Unidos
$ws.Cells.Item(34,2).Value = 53045
$ws.Cells.Item(34,3).Value = 445
$ws.Cells.Item(34,4).Value = 42282
$ws.Cells.Item(34,5).Value = 10436
$ws.Cells.Item(34,7).Value = 1
$ws.Cells.Item(34,8).Value = 327

# Row 49 - Suiza
$ws.Cells.Item(49,2).Value = 32498
$ws.Cells.Item(49,3).Value = 129
$ws.Cells.Item(49,4).Value = 29400
$ws.Cells.Item(49,5).Value = 1132

# Row 64 - Nepal
$ws.Cells.Item(64,2).Value = 16423
$ws.Cells.Item(64,3).Value = 255
$ws.Cells.Item(64,4).Value = 7752
$ws.Cells.Item(64,5).Value = 8636

# Row 82 - Finlandia
$ws.Cells.Item(82,4).Value = 6800
$ws.Cells.Item(82,5).Value = 136

# Rows 87-90 - Bosnia y Herzegovina overtakes Etiopia / Gabon / Guinea in ranking,
# so the table (sorted by total cases descending) re-shuffles these four countries.
# Row 87 becomes Bosnia y Herzegovina with freshly updated figures.
$ws.Cells.Item(87,1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(87,2).Value = 5869
$ws.Cells.Item(87,3).Value = 248
$ws.Cells.Item(87,4).Value = 2769
$ws.Cells.Item(87,5).Value = 2891
$ws.Cells.Item(87,7).Value = 2
$ws.Cells.Item(87,8).Value = 209

# Row 88 becomes Etiopia (previous row 87 figures)
$ws.Cells.Item(88,1).Value = "Etiopia"
$ws.Cells.Item(88,2).Value = 5846
$ws.Cells.Item(88,4).Value = 2430
$ws.Cells.Item(88,5).Value = 3313
$ws.Cells.Item(88,8).Value = 103

# Row 89 becomes Gabon (previous row 88 figures)
$ws.Cells.Item(89,1).Value = "Gabon"
$ws.Cells.Item(89,2).Value = 5743
$ws.Cells.Item(89,4).Value = 2574
$ws.Cells.Item(89,5).Value = 3123
$ws.Cells.Item(89,8).Value = 46

# Row 90 becomes Guinea with freshly updated figures
$ws.Cells.Item(90,1).Value = "Guinea"
$ws.Cells.Item(90,2).Value = 5697
$ws.Cells.Item(90,3).Value = 61
$ws.Cells.Item(90,4).Value = 4577
$ws.Cells.Item(90,5).Value = 1086
$ws.Cells.Item(90,8).Value = 34

# Row 100 - Madagascar
$ws.Cells.Item(100,2).Value = 3573
$ws.Cells.Item(100,3).Value = 101
$ws.Cells.Item(100,4).Value = 1761
$ws.Cells.Item(100,5).Value = 1779

# Row 117 - Islandia
$ws.Cells.Item(117,2).Value = 1880
$ws.Cells.Item(117,3).Value = 7
$ws.Cells.Item(117,4).Value = 1850
$ws.Cells.Item(117,5).Value = 20

# Row 153 - Malta
$ws.Cells.Item(153,4).Value = 654
$ws.Cells.Item(153,5).Value = 10

# Rows 209-210 - Groenlandia / Islas Malvinas swap position (tie on total cases,
# Groenlandia now listed first)
$ws.Cells.Item(209,1).Value = "Groenlandia"
$ws.Cells.Item(210,1).Value = "Islas Malvinas"
